# Daily attendance processing - 2026-01-02 19:52:51
# Normalizes the "Recorded By" (column G) values so the user's email
# address is listed before the literal "System" entry, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$updated = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $updated++
    }
}

Write-Host "Updated $updated 'Recorded By' cell(s) in column G."
